# Insert a new weekly price record as row 295 in the Zanahoria (carrot)
# price sheet, pushing all existing rows from 295 down to 296 (the prior
# last row, 416, becomes row 417). This mirrors the commit's weekly
# fruit/vegetable data append ("Fruta / hortaliza, semanal").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 295:416 down by one row, creating a blank row 295.
$ws.Rows.Item(295).Insert()

# Populate the newly inserted row 295 with the new record's data.
$ws.Cells.Item(295, 1).Value = 4
$ws.Cells.Item(295, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(295, 3).Value = "Los Lagos"
$ws.Cells.Item(295, 4).Value = 44784
$ws.Cells.Item(295, 5).Value = 10
$ws.Cells.Item(295, 6).Value = 100114013
$ws.Cells.Item(295, 7).Value = "Zanahoria"
$ws.Cells.Item(295, 8).Value = "Sin especificar"
$ws.Cells.Item(295, 9).Value = "Primera"
$ws.Cells.Item(295, 10).Value = 300
$ws.Cells.Item(295, 11).Value = 10000
$ws.Cells.Item(295, 12).Value = 10500
$ws.Cells.Item(295, 13).Value = 10250
$ws.Cells.Item(295, 14).Value = "`$/saco 20 kilos"
$ws.Cells.Item(295, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(295, 16).Value = 512
$ws.Cells.Item(295, 17).Value = 20
$ws.Cells.Item(295, 18).Value = "Hortaliza"
